$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: append a short run of text (punctuation) immediately before
# the end-of-paragraph mark of paragraph number $paraIndex, using the
# body-text ("normal list item") character formatting: Arial, not
# bold, size 24 half-points (12pt).  Nudging Font.Size away and back
# forces the engine to keep this as its own <w:r> instead of silently
# merging it back into the preceding run (which would happen if the
# resulting formatting were left byte-identical without ever being
# touched).
# ---------------------------------------------------------------------
function Add-BodyPunct($paraIndex, $text) {
    $para = $d.Paragraphs($paraIndex).Range
    $insertPos = $para.End - 1
    $ins = $d.Range($insertPos, $insertPos)
    $ins.InsertAfter($text)
    $newRange = $d.Range($insertPos, $insertPos + $text.Length)
    $newRange.Font.Size = 30
    $newRange.Font.Size = 12
}

# ---------------------------------------------------------------------
# Helper: append a short run of text (punctuation) immediately before
# the end-of-paragraph mark of paragraph number $paraIndex, using the
# heading character formatting: Arial, bold, size 28 half-points
# (14pt), szCs 24.
# ---------------------------------------------------------------------
function Add-HeadingPunct($paraIndex, $text) {
    $para = $d.Paragraphs($paraIndex).Range
    $insertPos = $para.End - 1
    $ins = $d.Range($insertPos, $insertPos)
    $ins.InsertAfter($text)
    $newRange = $d.Range($insertPos, $insertPos + $text.Length)
    $newRange.Font.Size = 30
    $newRange.Font.Size = 14
}

# 1) "LISTA DE POSSÍVEIS CENÁRIOS" heading -> add trailing "."
Add-HeadingPunct 30 "."

# 2) "Leitura de produtos" -> add trailing ";"
Add-BodyPunct 31 ";"

# 3) "Comunicar com banco de dados" -> add trailing ";"
Add-BodyPunct 32 ";"

# 4) "Listar produtos" -> add trailing ";"
Add-BodyPunct 33 ";"

# 5) "Gerar arquivos PDF" -> add trailing ";"
Add-BodyPunct 34 ";"

# 6) "Iniciar" + " impressão de arquivos" -> add trailing ";" (the two
#    existing runs are left as-is, a new run is appended)
Add-BodyPunct 35 ";"

# 7) "Criar " + "e controlar " + "histórico" + " de saídas" -> collapse
#    into a single run (re-assert the same text via Find/Replace, which
#    merges runs sharing identical formatting), then append ";"
$r36 = $d.Paragraphs(36).Range
$r36.Find.Execute("Criar e controlar histórico de saídas", $true, $false, $false, $false, $false, $true, 1, $false, "Criar e controlar histórico de saídas", 2) | Out-Null
Add-BodyPunct 36 ";"

# 8) "Limpar campos" -> add trailing ";"
Add-BodyPunct 37 ";"

# 9) "Leitura de produtos sem o leitor de códigos de barras" -> add trailing ";"
Add-BodyPunct 38 ";"

# 10) "Salvar arquivos PDF sem a presença de impressoras" -> add trailing "."
Add-BodyPunct 39 "."

# 11) "DESCRIÇÃO DE CADA CENÁRIO" heading -> add trailing "."
Add-HeadingPunct 41 "."

# 12) "Iniciar impressão de arquivos" + "." -> collapse into a single run
$r50 = $d.Paragraphs(50).Range
$r50.Find.Execute("Iniciar impressão de arquivos.", $true, $false, $false, $false, $false, $true, 1, $false, "Iniciar impressão de arquivos.", 2) | Out-Null

Write-Host "done"
